$wb = $excel.ActiveWorkbook

# --- Update template version on the "isa_template" sheet ---
$wsTemplate = $wb.Worksheets.Item("isa_template")
$wsTemplate.Range("B4").Value = "1.0.4"

# --- Update MIAPPE ontology term columns (0076 -> 0079) on the table sheet ---
$wsTable = $wb.Worksheets.Item("New Table")

$wsTable.Range("F1").Value = "Characteristic [Sample description]"
$wsTable.Range("G1").Value = "Term Source REF (MIAPPE:0079)"
$wsTable.Range("H1").Value = "Term Accession Number (MIAPPE:0079)"

# --- Update the underlying ListObject (table) column names to match ---
$table = $wsTable.ListObjects.Item("annotationTable")
$table.ListColumns.Item(6).Name = "Characteristic [Sample description]"
$table.ListColumns.Item(7).Name = "Term Source REF (MIAPPE:0079)"
$table.ListColumns.Item(8).Name = "Term Accession Number (MIAPPE:0079)"
